$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "kmeans_tuning_1"
$ws.Range("B4").Value = 0.5050436415844799
$ws.Range("C4").Value = 714477.0893068067
$ws.Range("D4").Value = 0.4315535954435768

# Match the formatting used on the other "Iteration name" column cells (A2, A3)
# by copying A3's format (bold, centered, top-aligned, thin border) onto A4.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
